$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update director names (titles/teams stay the same)
$ws.Range("A2").Value = "Yeseo Han"

# Row 3 (Maxwell Xu / Marketing Intern) is cleared out entirely
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""

$ws.Range("A4").Value = "Forest Huang"
$ws.Range("A5").Value = "Andrew Doan"
$ws.Range("A6").Value = "Jessie Yang"
$ws.Range("A7").Value = "David Ayala"
$ws.Range("A8").Value = "Julia Lin"
$ws.Range("A9").Value = "Jason Henkel"
$ws.Range("A10").Value = "Nathan Lee"
$ws.Range("A11").Value = "Winnie Qi"

# Move the active selection like in the saved file
$ws.Range("A3").Select() | Out-Null
